$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.127.93'
$ws.Range("E2").Value = '  +0.52%  '
$ws.Range("D3").Value = '1.788.68'
$ws.Range("E3").Value = '  +0.60%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '226.46'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.46%  '
$ws.Range("E6").Value = '  -0.71%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E8").Value = '  -1.26%  '
$ws.Range("E9").Value = '  +1.14%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0689'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.95%  '
$ws.Range("E11").Value = '  +0.88%  '
$ws.Range("D12").Value = '2.045.56'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.15'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.94%  '
$ws.Range("D14").Value = '1.776.86'
$ws.Range("E14").Value = '  -0.04%  '
$ws.Range("D15").Value = '34.061.53'
$ws.Range("E15").Value = '  +0.40%  '
$ws.Range("E16").Value = '  +0.01%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.18'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.98%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '68.01'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.30%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '245.34'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.15%  '
$ws.Range("D20").Value = '0.0₃0778'
$ws.Range("E20").Value = '  -0.66%  '
$ws.Range("E21").Value = '  -0.05%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.83'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.46%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.09'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.36%  '
$ws.Range("E24").Value = '  -1.02%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '161.29'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.70%  '
$ws.Range("E26").Value = '  +1.04%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.30'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.23%  '
$ws.Range("E28").Value = '  +0.49%  '
$ws.Range("E29").Value = '  +0.03%  '
$ws.Range("E30").Value = '  -0.64%  '
$ws.Range("E31").Value = '  +1.36%  '
$ws.Range("E32").Value = '  +0.44%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.61'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.95%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.81'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.00%  '
$ws.Range("D35").Value = '1.454.39'
$ws.Range("E35").Value = '  +4.57%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.43'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +9.87%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.646'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.18%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0191'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.49%  '
$ws.Range("E39").Value = '  -0.52%  '
$ws.Range("E40").Value = '  +3.35%  '
$ws.Range("E41").Value = '  +0.70%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.918'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.99%  '
$ws.Range("E43").Value = '  +0.40%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.52'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.43%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0510'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.81%  '
$ws.Range("E46").Value = '  +3.84%  '
$ws.Range("E47").Value = '  -0.18%  '
$ws.Range("E48").Value = '  -0.85%  '
$ws.Range("D49").Value = '1.947.23'
$ws.Range("E49").Value = '  +0.75%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '106.23'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.61%  '
$ws.Range("E51").Value = '  +0.02%  '
